# Add a new booking row (row 6) to the "All Bookings" sheet, matching the
# format/content of the existing booking rows (2-5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 6
$sourceRow = 5

$values = @(
    "2025-04-18T08-52-52-064Z",
    "kamal",
    "pavankumarbnm@gmail.com",
    "6729738922",
    "JSN Signature",
    "2025-04-30",
    "03:00 PM - 04:00 PM",
    "4/18/2025, 2:22:52 PM"
)

# 1) Copy the formatting (and values, temporarily) of the row above into the
#    new row, so the new row picks up the same style used by every other
#    booking data row.
$ws.Range("A$sourceRow`:H$sourceRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) Stage the real values in a scratch area far away from the used range,
#    formatted as Text, so that numeric-looking / date-looking strings
#    (e.g. the phone number and the visit date) are kept as plain text
#    instead of being auto-converted to numbers/dates by Excel.
$scratchRow = 500
$scratchRange = $ws.Range("A$scratchRow`:H$scratchRow")
$scratchRange.NumberFormat = "@"
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($scratchRow, $i + 1).Value = $values[$i]
}

# 3) Copy just the values (not the scratch formatting) from the staging row
#    into the new row, which already has the correct style applied.
$scratchRange.Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# 4) Clean up the scratch area.
$scratchRange.Clear()
